# Update the Alpha Xi leadership roster (names only; titles in column B are unchanged).
# Edit order matches the order the names were typed in the source workbook so that
# the shared-string table indices line up with the target file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sean Devine"
$ws.Range("A3").Value = "Jenny Wu"
$ws.Range("A7").Value = "Aryan Bajaria"
$ws.Range("A4").Value = "TBD"
$ws.Range("A5").Value = "Ronica Cheng"
$ws.Range("A6").Value = "Haley Truong"
$ws.Range("A8").Value = "Nathan Lee"
$ws.Range("A9").Value = "Madeline Li"
$ws.Range("A10").Value = "Jorina Chen"
$ws.Range("A11").Value = "Yeseo Han"

# Leave the final selection where the author left it when saving.
$ws.Range("D13").Select() | Out-Null
